# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (blue "Office" palette) - currently unused by
#                             any slide, only wired to the notes master.
#   ppt/theme/theme2.xml  -> "Integral" (red/violet palette) - the theme actually applied
#                             to the single slide master / all slides.
#
# The target edit swaps the two themes' content: the deck's live (slide-master) theme
# becomes the "Office Theme" palette, while the red/violet "Integral" palette moves over
# to the other theme part. The part that drives what slides actually look like is the one
# reachable from SlideMaster.ColorScheme, so recolor it to the Office Theme values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function ToComRgb($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), in the standard
# DrawingML clrScheme slot order that ColorScheme.Colors(1..12) maps onto.
$officeTheme = @(
    @(0x00, 0x00, 0x00),   # 1  dk1
    @(0xFF, 0xFF, 0xFF),   # 2  lt1
    @(0x44, 0x54, 0x6A),   # 3  dk2
    @(0xE7, 0xE6, 0xE6),   # 4  lt2
    @(0x5B, 0x9B, 0xD5),   # 5  accent1
    @(0xED, 0x7D, 0x31),   # 6  accent2
    @(0xA5, 0xA5, 0xA5),   # 7  accent3
    @(0xFF, 0xC0, 0x00),   # 8  accent4
    @(0x44, 0x72, 0xC4),   # 9  accent5
    @(0x70, 0xAD, 0x47),   # 10 accent6
    @(0x05, 0x63, 0xC1),   # 11 hlink
    @(0x95, 0x4F, 0x72)    # 12 folHlink
)

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $rgb = $officeTheme[$i]
    $colorScheme.Colors($i + 1).RGB = ToComRgb $rgb[0] $rgb[1] $rgb[2]
}
